$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update "last selected cell" bookkeeping on the sheets the author had
#    clicked around in before adding the new sheet. These are cosmetic
#    <selection> changes captured in the saved file.
# ---------------------------------------------------------------------------

$ws = $wb.Worksheets.Item("General")
$ws.Range("M41").Select()

$ws = $wb.Worksheets.Item("xdyt-bug-with")
$ws.Range("D67").Select()

$ws = $wb.Worksheets.Item("xdyt-bug-wo")
$ws.Range("H2:H7").Select()

$ws = $wb.Worksheets.Item("xdyt-impl-with")
$ws.Range("N2:N7").Select()

$ws = $wb.Worksheets.Item("xdyt-impl-wo")
$ws.Range("H2:H7").Select()

$ws = $wb.Worksheets.Item("xdc-impl-with")
$ws.Range("N2:N7").Select()

$ws = $wb.Worksheets.Item("xdc-impl-wo")
$ws.Range("H2:H7").Select()

$ws = $wb.Worksheets.Item("xdc-bug-with")
$ws.Range("N2:N7").Select()

$ws = $wb.Worksheets.Item("xdc-bug-wo")
$ws.Range("J48").Select()

$ws = $wb.Worksheets.Item("Debugging strategies")
$ws.Range("I57").Select()

# ---------------------------------------------------------------------------
# 2. Add the new "number of devices used" worksheet as the last tab.
# ---------------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$new.Name = "number of devices used"

# Row 1
$new.Range("A1").Value = "Total"
$new.Range("K1").Value = "Emulated"

# Row 2 headers -- write order matters: it drives the order new shared
# strings are appended in (xdcb-wo, xdci-w, xdci-wo, xdcb-w, xdytb-w,
# xdytb-wo, xdyti-w, xdyti-wo).
$new.Range("C2").Value = "xdcb-wo"
$new.Range("D2").Value = "xdci-w"
$new.Range("E2").Value = "xdci-wo"
$new.Range("B2").Value = "xdcb-w"
$new.Range("F2").Value = "xdytb-w"
$new.Range("G2").Value = "xdytb-wo"
$new.Range("H2").Value = "xdyti-w"
$new.Range("I2").Value = "xdyti-wo"

$new.Range("L2").Value = "xdcb-w"
$new.Range("M2").Value = "xdcb-wo"
$new.Range("N2").Value = "xdci-w"
$new.Range("O2").Value = "xdci-wo"
$new.Range("P2").Value = "xdytb-w"
$new.Range("Q2").Value = "xdytb-wo"
$new.Range("R2").Value = "xdyti-w"
$new.Range("S2").Value = "xdyti-wo"

# Rows 3-8: raw data (two copies, B:I = "real" devices, L:S = "emulated")
$dataReal = @(
    @(2,3,3,3,5,2,2,2),
    @(1,1,3,3,2,3,3,2),
    @(2,2,3,3,3,2,3,2),
    @(2,1,3,3,2,3,4,3),
    @(1,1,3,3,2,2,2,2),
    @(3,1,3,1,2,2,2,3)
)
$dataEmulated = @(
    @(2,3,3,3,4,2,1,2),
    @(1,1,3,3,2,3,3,2),
    @(2,2,3,3,3,2,3,2),
    @(2,1,3,3,2,3,4,3),
    @(1,1,3,3,2,2,2,2),
    @(3,1,3,1,2,1,2,2)
)
$colsReal = @("B","C","D","E","F","G","H","I")
$colsEmulated = @("L","M","N","O","P","Q","R","S")

for ($r = 0; $r -lt 6; $r++) {
    $row = 3 + $r
    for ($c = 0; $c -lt 8; $c++) {
        $new.Range($colsReal[$c] + $row).Value = $dataReal[$r][$c]
        $new.Range($colsEmulated[$c] + $row).Value = $dataEmulated[$r][$c]
    }
}

# Row 12: Avg
$new.Range("A12").Value = "Avg"
for ($c = 0; $c -lt 8; $c++) {
    $col = $colsReal[$c]
    $new.Range($col + "12").Formula = "=AVERAGE(" + $col + "3:" + $col + "8)"
}
for ($c = 0; $c -lt 8; $c++) {
    $col = $colsEmulated[$c]
    $new.Range($col + "12").Formula = "=AVERAGE(" + $col + "3:" + $col + "8)"
}

# Row 13: StDev
$new.Range("A13").Value = "StDev"
for ($c = 0; $c -lt 8; $c++) {
    $col = $colsReal[$c]
    $new.Range($col + "13").Formula = "=STDEV(" + $col + "3:" + $col + "8)"
}
for ($c = 0; $c -lt 8; $c++) {
    $col = $colsEmulated[$c]
    $new.Range($col + "13").Formula = "=STDEV(" + $col + "3:" + $col + "8)"
}

$new.Range("O21").Select()
$new.Activate()
